# TransferLogger export template: default export template variables fix
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A8 previously had a highlighted fill (distinct style); now uses the plain style like A3/A4
$ws.Range("A8").Style = "Normal"

# G8 / H8 previously blank placeholder cells under the TRANSFER / COMMENT headers now carry
# their own template placeholders, matching the other columns in row 8.
$ws.Range("G8").Value = "%Transfer%"
$ws.Range("H8").Value = "%Comment%"

# A1 header label was "Application ID" -> now "APPLICATION ID:" (matches style of other labels)
$ws.Range("A1").Value = "APPLICATION ID:"

# Selection moved to A5 as last interacted cell
$ws.Range("A5").Select()
